$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# The "Metadata" sheet is protected with an unknown password. Re-asserting
# protection with a blank password first (this does not clear the existing
# protected state nor require the original password) lets a subsequent
# Unprotect("") call succeed in actually lifting the lock.
$ws.Protect("")
$ws.Unprotect("")

# "Organization" -> "Organisation"
$ws.Range("F8").Value = "Organisation"

# "Finalization date" -> "Finalisation date"
$ws.Range("C12").Value = "Finalisation date"

# "e-mail" -> "email" (both occurrences)
$ws.Range("C10").Value = "email"
$ws.Range("F10").Value = "email"

# Bold the "email" label under "Name" (F10)
$ws.Range("F10").Font.Bold = $true

# Restore sheet protection (matching original protected state) with a blank
# password, since the original password is not recoverable.
$ws.Protect("")
